# Rename "cut_before"/"cut_after" option names to "cut_ends"/"cut_beginnings"
# throughout the regression-rules example sheet, and update the selected
# cell on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the explanatory text cells (B8, B9, B10) and the column heading
# (O13) that referenced the old option names.
$ws.Range("B8").Value  = "cut_beginnings is the number seconds of data which are cut out after a switch between chambers"
$ws.Range("B9").Value  = "cut_ends is the number of seconds of data which are cut out before a switch between chambers"
$ws.Range("B10").Value = "cut_beginnings and cut_ends may not be side or substance or side specific (so N2O:cut_ends or left:cut_ends is not allowed)"
$ws.Range("O13").Value = "cut_beginnings"

# Move the selection to match the edited workbook (the author ended up with
# O13 selected instead of the old B17:B19 block).
$ws.Activate() | Out-Null
$ws.Range("O13").Select() | Out-Null
